$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells for team record columns (AD, AE, AF) on row 1,
# reusing the existing header style (bold/centered/bordered) from AC1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill the team record (Wins/Losses/Ties) for every data row (2-43).
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 30).Value = 85
    $ws.Cells.Item($r, 31).Value = 77
    $ws.Cells.Item($r, 32).Value = 0
}
